$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.449.57'
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('D3').Value = '1.877.50'
$ws.Range('E3').Value = '  +0.96%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7128'
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.97'
$ws.Range('E6').Value = '  +1.51%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07840'
$ws.Range('E8').Value = '  -2.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3113'
$ws.Range('E9').Value = '  +2.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.16'
$ws.Range('E10').Value = '  +6.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08257'
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7302'
$ws.Range('E12').Value = '  +2.96%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.867.89'
$ws.Range('E13').Value = '  -3.69%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.269'
$ws.Range('E14').Value = '  +1.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.09'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '29.455.34'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.916'
$ws.Range('E17').Value = '  +1.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '247.60'
$ws.Range('E18').Value = '  +3.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007871'
$ws.Range('E19').Value = '  -0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.28'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.951'
$ws.Range('E22').Value = '  +6.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.1581'
$ws.Range('E24').Value = '  +9.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '163.70'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.003'
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.29'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.362'
$ws.Range('E28').Value = '  -3.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.494'
$ws.Range('E29').Value = '  +1.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.376'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.126'
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05318'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.932'
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.200'
$ws.Range('E34').Value = '  +3.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7231'
$ws.Range('E35').Value = '  +0.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.679'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01865'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').Value = '1.256.74'
$ws.Range('E38').Value = '  +8.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.736'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9095'
$ws.Range('E40').Value = '  -3.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.93'
$ws.Range('E41').Value = '  +4.09%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.142'
$ws.Range('E42').Value = '  +2.38%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.48'
$ws.Range('E44').Value = '  +0.47%  '
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('B46').Value = 'SynthetixNetwork'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.949'
$ws.Range('E46').Value = '  +13.10%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.768'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('B48').Value = 'BabyDogeCoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.00000000120'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('B49').Value = 'TheSandbox'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4323'
$ws.Range('E49').Value = '  +1.09%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.235'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.076'
$ws.Range('E51').Value = '  +1.36%  '
